$d = $word.ActiveDocument

# --- Edit 1: first paragraph gets two trailing spaces, then three new
#     red-colored runs appended forming "(This is a change – Version for main branch)" ---
$p1 = $d.Paragraphs.First
$r1 = $p1.Range
$r1.End = $r1.End - 1          # exclude the paragraph mark so inserts stay in paragraph 1
$r1.Collapse(0)
$r1.InsertAfter("  ")

$r1.Collapse(0)
$r1.InsertAfter("(This is a change " + [char]0x2013 + " Ve")
$r1.Font.Color = 255

$r1.Collapse(0)
$r1.InsertAfter("rsion for main branch")
$r1.Font.Color = 255

$r1.Collapse(0)
$r1.InsertAfter(")")
$r1.Font.Color = 255

# --- Edit 2: remove the trailing paragraph "ank God almighty, we are free at last." ---
$last = $d.Paragraphs.Last
$last.Range.Delete()
